$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2021-04-22", "overview", "K02000001", "United Kingdom", 4398431, 2729, 18, 127345),
    @("2021-04-23", "overview", "K02000001", "United Kingdom", 4401109, 2678, 40, 127385),
    @("2021-04-24", "overview", "K02000001", "United Kingdom", 4403170, 2061, 32, 127417),
    @("2021-04-25", "overview", "K02000001", "United Kingdom", 4404882, 1712, 11, 127428),
    @("2021-04-26", "overview", "K02000001", "United Kingdom", 4406946, 2064, 6, 127434),
    @("2021-04-27", "overview", "K02000001", "United Kingdom", 4409631, 2685, 17, 127451)
)

$startRow = 254
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = "'" + $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
}
